$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "2020" column (H) mirroring column G's formatting ---
# Copy formatting from column G (rows 4-25) into column H first so the
# new cells inherit the same number formats / borders / fonts as the
# existing year columns.
$ws.Range("G4:G25").Copy()
$ws.Range("H4:H25").PasteSpecial(-4122)

# Rows 8 and 20 use the "0.0" decimal style (same style already used by
# G10/G22/G23) rather than the plain style used elsewhere in column G.
$ws.Range("G10").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("H20").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Populate the new column H values ---
$ws.Range("H4").Value = 2020
$ws.Range("H5").Value = 42.2
$ws.Range("H7").Value = 42.5
$ws.Range("H8").Value = 42
$ws.Range("H10").Value = 50.9
$ws.Range("H11").Value = 36.9
$ws.Range("H12").Value = 34.799999999999997
$ws.Range("H14").Value = 30.7
$ws.Range("H15").Value = 48.8
$ws.Range("H17").Value = 61.1
$ws.Range("H18").Value = 56.7
$ws.Range("H19").Value = 41.6
$ws.Range("H20").Value = 49
$ws.Range("H21").Value = 43.5
$ws.Range("H22").Value = 33.9
$ws.Range("H23").Value = 34.6
$ws.Range("H24").Value = 23.6
$ws.Range("H25").Value = 35.9

# --- Update the sheet view: drop the frozen topLeftCell offset and move
#     the active selection from H15 to B13 ---
$ws.Range("B13").Select()
